$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O2").Value = "-- se calcoliamo un anno solo forecast: mettere mese_fine_cons e max_emissione_cons a 0"
$ws.Range("O3").Value = "-- se non vogliamo calcolare rateo anni precedenti, inserire 00000000 sotto ds_ric_gas_anno_prec e ds_ric_pwr_anno_prec"
$ws.Range("O4").Value = "-- anno_prec_tipologia = DELTA o anno_prec_tipologia = EMISSIONE -> modalità di calcolo del rateo anni precedenti (utile nel BUDGET fatto a ottobre/novembre sull'anno successivo) "
$ws.Range("O1").Value = "GUIDA"
$ws.Range("O1").Font.Bold = $true
$ws.Range("O1").Select() | Out-Null
